$d = $word.ActiveDocument

# --- Change 1: remove spell-check proofErr wrapper runs around
# "Shneiderman" and "Additon" by collapsing the whole citation sentence
# (currently split across five runs, two wrapped in <w:proofErr/>) into a
# single run with identical text. Replacing a Range spanning multiple
# runs in one shot merges them into one run and drops the now-orphaned
# <w:proofErr/> markers. (The engine no-ops a Range.Text assignment whose
# value already equals the range's current text, so first stamp a throwaway
# value to force the merge, then set the real text.)
$citationPara = $d.Paragraphs.Item(3)
$citationRange = $d.Range($citationPara.Range.Start, $citationPara.Range.End - 1)
$citationRange.Text = "placeholder"
$citationPara = $d.Paragraphs.Item(3)
$citationRange = $d.Range($citationPara.Range.Start, $citationPara.Range.End - 1)
$citationRange.Text = "Cognitive walkthrough, heuristic evaluation, review based see Shneiderman, B. (1998) Designing the user interface: Strategies for effective human computer interaction (3rd ed.). Reading, MA: Additon-Wesley Publishing"

# --- Change 2: drop the _GoBack bookmark from the now-empty paragraph
# right after the citation (it will be re-added further down, after the
# newly-written paragraph, per the diff).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Change 3: add a new paragraph after "Observed how clusters tended
# to flatten the data." with the walkthrough note, and re-create the
# _GoBack bookmark (collapsed, at the end of that new paragraph's text).
$clusterPara = $d.Paragraphs.Item(18)
$clusterPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(19)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$noteText = "Zoomed in on one or two locations on the map where by experience certain strategies should be useful, if those strategies weren" + [char]0x2019 + "t linked to that location the results were dismissed as incorrect "
# Append a throwaway trailing marker character so the bookmark's insertion
# point below isn't exactly on the paragraph-end boundary (collapsed
# ranges anchored right before a paragraph mark misplace on this engine);
# the marker is stripped immediately after the bookmark is anchored.
$newRange.Text = $noteText + "Z"

$newPara = $d.Paragraphs.Item(19)
$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$markerRange.Delete()
